# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2310"
#   "<header>_new" -> "<header>_FV2404"
# Then freeze the header row and wrap the data range in an Excel Table
# (ListObject) so the new headers double as the table's column headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base (suffix-less) header names, in left-to-right column order, shared by
# both the "_FV2310" (old) and "_FV2404" (new) blocks of columns.
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J hold the "_old" (-> "_FV2310") headers, columns L-U hold the
# "_new" (-> "_FV2404") headers; column K ("diff") is left untouched.
$fv2310Cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$fv2404Cols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($fv2310Cols[$i] + "1").Value = $baseNames[$i] + "_FV2310"
    $ws.Range($fv2404Cols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
}

# Freeze the header row (row 1) in the sheet view.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Wrap the used range in a Table so the headers become the table's column
# headers (ListObject picks up the renamed cells from row 1 automatically).
$dataRange = $ws.Range("A1:U58")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)

Write-Output "Renamed headers, froze top row, and created table $($table.Name)"
